# Apply the update described in the commit: add a new "CS30-RT" item row
# to the "Instal days by Model" table, right after the existing "CS30" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instal days by Model")

# The existing table ("Table1") covers A1:E41 (header + 40 data rows).
# "CS30" currently lives on row 19; insert a new row above it so the new
# "CS30-RT" entry lands at row 19 and everything else shifts down by one.
$ws.Rows.Item(19).Insert()

$ws.Range("A19").Value = "CS30-RT"
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = $true
$ws.Range("E19").Value = $true

# Grow the table / autofilter range to include the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E42"))

# Keep the workbook-level filter-database defined name in sync with the
# new table extent.
$fd = $wb.Names.Item("_xlnm._FilterDatabase")
$fd.RefersTo = "='Instal days by Model'!`$A`$1:`$C`$42"

# Reflect the selection left behind after the edit.
$ws.Activate()
$ws.Range("C20").Select()
